$wb = $excel.ActiveWorkbook

# "Repayment schedule" sheet (sheet3.xml) gains a new blank column before the
# existing "Late" column (old column N), shifting Late/Outstanding(heading)/Outstanding
# one column to the right (N->O, O->P, P->Q), and becomes the active tab with a
# new selection. "NewLoanInput" (sheet1.xml) loses the tabSelected flag as the
# active sheet moves away from it.
$wsSchedule = $wb.Worksheets.Item("Repayment schedule")

# Insert a new column at N, pushing existing N:P data to O:Q.
$wsSchedule.Columns("N:N").Insert()

# The freshly inserted column picks up the same width as column M (11 chars-ish),
# matching what Excel does when inserting a column (inherits left neighbour's width).
$wsSchedule.Columns("N:N").ColumnWidth = $wsSchedule.Columns("M:M").ColumnWidth

# Make "Repayment schedule" the active sheet/tab, with the new selection.
$wsSchedule.Activate()
$wsSchedule.Range("K18").Select()
